$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename three existing category labels (column A) ---
$ws.Range("A7").Value  = "Arterial hypertension"   # was "Hypertoni"
$ws.Range("A33").Value = "CNS disease"              # was "CNS"
$ws.Range("A39").Value = "Rheumatic disease "       # was "Reuma" (trailing space kept)

# --- Append five new ECI categories as rows 45-49 ---
$ws.Range("A45").Value = "Hypothyroidism"
$ws.Range("B45").Value = "ECI"
$ws.Range("C45").Value = "hypothyroidism"

$ws.Range("A46").Value = "Coagulopathy"
$ws.Range("B46").Value = "ECI"
$ws.Range("C46").Value = "coagulopathy"

$ws.Range("A47").Value = "Obesity"
$ws.Range("B47").Value = "ECI"
$ws.Range("C47").Value = "obesity"

$ws.Range("A48").Value = "Weight loss"
$ws.Range("B48").Value = "ECI"
$ws.Range("C48").Value = "weight loss"

$ws.Range("A49").Value = "Fluid electrolyte disorders"
$ws.Range("B49").Value = "ECI"
$ws.Range("C49").Value = "fluid electrolyte disorders"

# --- Widen column A so the longer labels keep fitting ---
$ws.Columns.Item(1).ColumnWidth = 24.7109375

# --- Page setup: A4, portrait (mirrors the Page Setup dialog being touched) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Move the selection to where editing left off ---
$ws.Range("A43").Select()
